$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights (also flips customHeight="1" like real Excel does) ---
$ws.Rows.Item(1).RowHeight = 29.25
$ws.Rows.Item(4).RowHeight = 26.25
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 28.5
$ws.Rows.Item(7).RowHeight = 41.25

# --- New column N, row 3: year header 2023 (same style as M3) ---
$ws.Range("N3").Value = 2023
$ws.Range("N3").Style = $ws.Range("M3").Style

# --- Row 4: M4 revised, N4 added (same style as M4) ---
$ws.Range("M4").Value = 923.8
$ws.Range("N4").Value = 583.5
$ws.Range("N4").Style = $ws.Range("M4").Style

# --- Row 5: N5 added empty, matching M5 style ---
$ws.Range("N5").Style = $ws.Range("M5").Style

# --- Row 6: N6 added empty, matching M6 style ---
$ws.Range("N6").Style = $ws.Range("M6").Style

# --- Row 7: M7 revised, N7 added (same style as M7) ---
$ws.Range("M7").Value = 64.03
$ws.Range("N7").Value = 64.08
$ws.Range("N7").Style = $ws.Range("M7").Style

# --- Reset the saved cursor/selection to the default top-left cell so the
#     stale "M14" selection isn't persisted on save ---
$ws.Range("A1").Select()
